$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.01136233333333333
$ws.Range("H2").Value = 0.034087
$ws.Range("I2").Value = 0.03319472577126831
$ws.Range("J2").Value = 0.03319472577126831
$ws.Range("N2").Value = 3.493359
$ws.Range("O2").Value = 0.9853654560111834
$ws.Range("P2").Value = 0.9853654560111834
$ws.Range("Q2").Value = 0.013230903137
$ws.Range("R2").Value = 0.119078128233
$ws.Range("S2").Value = 0.03270893609677198
$ws.Range("T2").Value = 0.03270893609677198

# Row 3
$ws.Range("G3").Value = 0.01136233333333333
$ws.Range("H3").Value = 0.034087
$ws.Range("I3").Value = 0.03319472577126831
$ws.Range("J3").Value = 0.03319472577126831
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01729433333333333
$ws.Range("N3").Value = 0.051883
$ws.Range("O3").Value = 0.01463454398881656
$ws.Range("P3").Value = 0.01463454398881656
$ws.Range("Q3").Value = 0.0001965039801111111
$ws.Range("R3").Value = 0.001768535821
$ws.Range("S3").Value = 0.0004857896744963287
$ws.Range("T3").Value = 0.0004857896744963287

# Row 4
$ws.Range("G4").Value = 0.330931
$ws.Range("H4").Value = 0.992793
$ws.Range("I4").Value = 0.9668052742287317
$ws.Range("J4").Value = 0.9668052742287317
$ws.Range("N4").Value = 3.493359
$ws.Range("O4").Value = 0.9853654560111834
$ws.Range("P4").Value = 0.9853654560111834
$ws.Range("Q4").Value = 0.385353595743
$ws.Range("R4").Value = 3.468182361687
$ws.Range("S4").Value = 0.9526565199144115
$ws.Range("T4").Value = 0.9526565199144115

# Row 5
$ws.Range("G5").Value = 0.330931
$ws.Range("H5").Value = 0.992793
$ws.Range("I5").Value = 0.9668052742287317
$ws.Range("J5").Value = 0.9668052742287317
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.01729433333333333
$ws.Range("N5").Value = 0.051883
$ws.Range("O5").Value = 0.01463454398881656
$ws.Range("P5").Value = 0.01463454398881656
$ws.Range("Q5").Value = 0.005723231024333334
$ws.Range("R5").Value = 0.051509079219
$ws.Range("S5").Value = 0.01414875431432023
$ws.Range("T5").Value = 0.01414875431432023
